# "fix: empty cell missing bug" — the sheet carried a block of trailing
# empty cells/rows (columns D:G on rows 1-4, and rows 5-21 entirely) that
# only existed to hold style references, with no real data. Trim the
# worksheet back down to its real data extent (A1:C4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 throwaway columns at the far right edge first so that, once we
# delete the 4 empty columns (D:G) below, the sheet's trailing default
# column-range ends up back at the full 16384 columns (min=4,max=16384)
# instead of shrinking to 16380 — this keeps the <cols> element identical
# in shape/width to the original, just re-split at column C/D instead of
# column G/H.
$ws.Columns("XFA:XFD").Insert()

# Drop the empty D:G columns (only ever held blank, style-only cells).
$ws.Columns("D:G").Delete()

# Drop the empty trailing rows 5:21 (only ever held blank, style-only
# cells) so the used range collapses to A1:C4.
$ws.Range("A5:G21").EntireRow.Delete()
